$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.059843338016958
$ws.Range("D2").Value = 1.057488595582969
$ws.Range("E2").Value = 1.064954358630998
$ws.Range("F2").Value = 1.074008256550537
$ws.Range("I2").Value = 1.051316569052303
$ws.Range("J2").Value = 1.064827267394267
$ws.Range("K2").Value = 1.060223427134904
$ws.Range("L2").Value = 1.067668902807052
$ws.Range("M2").Value = 1.076698606907402
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.061114870383562
$ws.Range("D3").Value = 1.05845684050984
$ws.Range("E3").Value = 1.06612946171137
$ws.Range("F3").Value = 1.075354904110726
$ws.Range("I3").Value = 1.051749249178972
$ws.Range("J3").Value = 1.065750676083209
$ws.Range("K3").Value = 1.061005126035141
$ws.Range("L3").Value = 1.068658410990062
$ws.Range("M3").Value = 1.077860996936937
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.06193697830824
$ws.Range("D4").Value = 1.059082677099288
$ws.Range("E4").Value = 1.066889439369922
$ws.Range("F4").Value = 1.076226201130028
$ws.Range("I4").Value = 1.052027541051118
$ws.Range("J4").Value = 1.066347004910497
$ws.Range("K4").Value = 1.061509636208021
$ws.Range("L4").Value = 1.069297718326627
$ws.Range("M4").Value = 1.078612520720357
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.062282438033622
$ws.Range("D5").Value = 1.059345617109902
$ws.Range("E5").Value = 1.067208842431638
$ws.Range("F5").Value = 1.076592479655904
$ws.Range("I5").Value = 1.052144133647642
$ws.Range("J5").Value = 1.066597421605677
$ws.Range("K5").Value = 1.061721422496387
$ws.Range("L5").Value = 1.069566253022882
$ws.Range("M5").Value = 1.0789283154603
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.062340433278001
$ws.Range("D6").Value = 1.059389756463154
$ws.Range("E6").Value = 1.067262466313841
$ws.Range("F6").Value = 1.076653978654047
$ws.Range("I6").Value = 1.052163686558266
$ws.Range("J6").Value = 1.066639451317128
$ws.Range("K6").Value = 1.061756964230511
$ws.Range("L6").Value = 1.069611327751986
$ws.Range("M6").Value = 1.078981330310139
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.061941594960865
$ws.Range("D7").Value = 1.05908619114977
$ws.Range("E7").Value = 1.066893707605931
$ws.Range("F7").Value = 1.07623109541583
$ws.Range("I7").Value = 1.05202910054253
$ws.Range("J7").Value = 1.066350352087946
$ws.Range("K7").Value = 1.061512467321866
$ws.Range("L7").Value = 1.069301307401735
$ws.Range("M7").Value = 1.078616740954963
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.060273196622991
$ws.Range("D8").Value = 1.057815960874312
$ws.Range("E8").Value = 1.06535157298568
$ws.Range("F8").Value = 1.074463379512884
$ws.Range("I8").Value = 1.051463144105224
$ws.Range("J8").Value = 1.0651395829654
$ws.Range("K8").Value = 1.060487876428536
$ws.Range("L8").Value = 1.068003514067316
$ws.Range("M8").Value = 1.077091572086987
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.057328077904753
$ws.Range("D9").Value = 1.055572346752933
$ws.Range("E9").Value = 1.062631017901788
$ws.Range("F9").Value = 1.071347737967044
$ws.Range("I9").Value = 1.050452929421842
$ws.Range("J9").Value = 1.062996931572489
$ws.Range("K9").Value = 1.058672371640558
$ws.Range("L9").Value = 1.065709100971358
$ws.Range("M9").Value = 1.074399160397426
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.055360974121484
$ws.Range("D10").Value = 1.054072925017169
$ws.Range("E10").Value = 1.060815063024425
$ws.Range("F10").Value = 1.069269983939147
$ws.Range("I10").Value = 1.049770687611688
$ws.Range("J10").Value = 1.061562227159309
$ws.Range("K10").Value = 1.057455170018462
$ws.Range("L10").Value = 1.064174281820859
$ws.Range("M10").Value = 1.072600777576274
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.05450827505
$ws.Range("D11").Value = 1.053422760687493
$ws.Range("E11").Value = 1.06002816511083
$ws.Range("F11").Value = 1.068370094349854
$ws.Range("I11").Value = 1.049473173482858
$ws.Range("J11").Value = 1.060939465915978
$ws.Range("K11").Value = 1.056926454525801
$ws.Range("L11").Value = 1.063508419673575
$ws.Range("M11").Value = 1.07182120363285
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.054191400859504
$ws.Range("D12").Value = 1.053181122469866
$ws.Range("E12").Value = 1.059735786291482
$ws.Range("F12").Value = 1.068035800214284
$ws.Range("I12").Value = 1.049362346626016
$ws.Range("J12").Value = 1.060707912804527
$ws.Range("K12").Value = 1.056729814501057
$ws.Range("L12").Value = 1.063260894317549
$ws.Range("M12").Value = 1.071531502242106
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.054259378020875
$ws.Range("D13").Value = 1.053232960980351
$ws.Range("E13").Value = 1.059798506659345
$ws.Range("F13").Value = 1.068107509088996
$ws.Range("I13").Value = 1.049386133726242
$ws.Range("J13").Value = 1.06075759226973
$ws.Range("K13").Value = 1.056772005873751
$ws.Range("L13").Value = 1.06331399819526
$ws.Range("M13").Value = 1.071593650242806
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.054482085087259
$ws.Range("D14").Value = 1.053402789630371
$ws.Range("E14").Value = 1.060003998838502
$ws.Range("F14").Value = 1.068342462229673
$ws.Range("I14").Value = 1.049464018973031
$ws.Range("J14").Value = 1.060920330409481
$ws.Range("K14").Value = 1.056910205344775
$ws.Range("L14").Value = 1.063487963135148
$ws.Range("M14").Value = 1.071797259556504
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.054619283086457
$ws.Range("D15").Value = 1.053507408280012
$ws.Range("E15").Value = 1.060130597345738
$ws.Range("F15").Value = 1.068487219853426
$ws.Range("I15").Value = 1.04951196461108
$ws.Range("J15").Value = 1.061020567941766
$ws.Range("K15").Value = 1.056995321198873
$ws.Range("L15").Value = 1.063595122815064
$ws.Range("M15").Value = 1.071922692260883
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.055417545008493
$ws.Range("D16").Value = 1.054116055023059
$ws.Range("E16").Value = 1.060867274413096
$ws.Range("F16").Value = 1.069329701935845
$ws.Range("I16").Value = 1.049790388292821
$ws.Range("J16").Value = 1.06160352541008
$ws.Range("K16").Value = 1.057490223973507
$ws.Range("L16").Value = 1.064218445781365
$ws.Range("M16").Value = 1.072652496876571
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.055918021743466
$ws.Range("D17").Value = 1.054497599179297
$ws.Range("E17").Value = 1.061329215657172
$ws.Range("F17").Value = 1.069858110242154
$ws.Range("I17").Value = 1.049964473097653
$ws.Range("J17").Value = 1.061968789047269
$ws.Range("K17").Value = 1.057800217598442
$ws.Range("L17").Value = 1.064609096322195
$ws.Range("M17").Value = 1.073110050792981
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.056209851884328
$ws.Range("D18").Value = 1.05472006028493
$ws.Range("E18").Value = 1.06159860266589
$ws.Range("F18").Value = 1.070166302017517
$ws.Range("I18").Value = 1.050065811407155
$ws.Range("J18").Value = 1.062181694218565
$ws.Range("K18").Value = 1.057980871705877
$ws.Range("L18").Value = 1.064836833320464
$ws.Range("M18").Value = 1.073376851279343
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.056309343366603
$ws.Range("D19").Value = 1.054795899060567
$ws.Range("E19").Value = 1.061690447354272
$ws.Range("F19").Value = 1.070271384194959
$ws.Range("I19").Value = 1.050100330865074
$ws.Range("J19").Value = 1.062254264541904
$ws.Range("K19").Value = 1.058042443029687
$ws.Range("L19").Value = 1.064914465000926
$ws.Range("M19").Value = 1.073467809378155
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.055864334611333
$ws.Range("D20").Value = 1.054456672133725
$ws.Range("E20").Value = 1.061279659498365
$ws.Range("F20").Value = 1.069801419113466
$ws.Range("I20").Value = 1.049945816393657
$ws.Range("J20").Value = 1.06192961493448
$ws.Range("K20").Value = 1.057766974775885
$ws.Range("L20").Value = 1.064567195945603
$ws.Range("M20").Value = 1.07306096819575
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.054416507408871
$ws.Range("D21").Value = 1.053352783160342
$ws.Range("E21").Value = 1.059943489020393
$ws.Range("F21").Value = 1.068273275384354
$ws.Range("I21").Value = 1.049441092468312
$ws.Range("J21").Value = 1.060872414506236
$ws.Range("K21").Value = 1.056869516007457
$ws.Range("L21").Value = 1.063436740189035
$ws.Range("M21").Value = 1.071737305384714
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.05350536705057
$ws.Range("D22").Value = 1.052657923405978
$ws.Range("E22").Value = 1.059102863924543
$ws.Range("F22").Value = 1.067312265963712
$ws.Range("I22").Value = 1.049121918476127
$ws.Range("J22").Value = 1.060206367797923
$ws.Range("K22").Value = 1.05630379153409
$ws.Range("L22").Value = 1.062724850983606
$ws.Range("M22").Value = 1.070904295343448
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.053988460194385
$ws.Range("D23").Value = 1.053026358398878
$ws.Range("E23").Value = 1.059548545693479
$ws.Range("F23").Value = 1.067821735914489
$ws.Range("I23").Value = 1.049291293002798
$ws.Range("J23").Value = 1.060559579975782
$ws.Range("K23").Value = 1.056603831606476
$ws.Range("L23").Value = 1.063102344700493
$ws.Range("M23").Value = 1.071345963906492
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.055888593791706
$ws.Range("D24").Value = 1.054475165573181
$ws.Range("E24").Value = 1.061302051964834
$ws.Range("F24").Value = 1.0698270354555
$ws.Range("I24").Value = 1.049954247180501
$ws.Range("J24").Value = 1.061947316483723
$ws.Range("K24").Value = 1.057781996270579
$ws.Range("L24").Value = 1.064586129301717
$ws.Range("M24").Value = 1.073083146764971
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.058090097352471
$ws.Range("D25").Value = 1.0561530150182
$ws.Range("E25").Value = 1.063334731977655
$ws.Range("F25").Value = 1.072153308375986
$ws.Range("I25").Value = 1.050715633650705
$ws.Range("J25").Value = 1.063551953368841
$ws.Range("K25").Value = 1.059142924953881
$ws.Range("L25").Value = 1.066303169535267
$ws.Range("M25").Value = 1.075095807636302
